# Updated symbol list on Mon Dec 12 22:47:31 UTC 2022 with GitHub Actions
#
# All "Price" (column D) values in this sheet are stored as text, not
# numbers (e.g. "0.06200" keeps a trailing zero). Assigning a plain
# numeric-looking string to .Value would make Excel coerce it into a
# real number and lose the exact textual formatting, so each numeric
# price value below is written with a leading apostrophe (Excel's
# "force text" prefix) and the cell style is then reset back to
# "Normal" so no stray number-format/quote-prefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'6.266"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06201"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.557"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.532"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'6.575"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8261"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1661"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08296"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03507"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03180"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09160"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.761"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001632"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04689"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.006279"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006222"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001069"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'3.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'0.01397"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'0.1243"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'0.0002736"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.04755"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1124"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003518"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.01140"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006333"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7229"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = "'0.00001900"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01240"
$ws.Range("D50").Style = "Normal"
